$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A135").Value = "Market Timing1"
$ws.Range("B135").Value = "Test market timing with HM model"
$ws.Range("C135").Value = "Market_Timing_test1"

$ws.Range("A136").Value = "Market Timing2"
$ws.Range("B136").Value = "Test market timing with TM model"
$ws.Range("C136").Value = "Market_Timing_test2"

$ws.Activate() | Out-Null
$ws.Range("B143").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 120
